# Auto-generated edit script: updates cryptos price/volume table
# to match the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.786.32"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "3.308.45"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'606.27"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").Value = "'141.66"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.307.12"
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").Value = "'5.55"
$ws.Range("E11").Value = "  +3.69%  "
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "'34.99"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "3.853.16"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").Value = "3.306.48"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").Value = "63.852.50"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("D20").Value = "'481.12"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").Value = "'8.00"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'13.92"
$ws.Range("E24").Value = "  +5.72%  "
$ws.Range("D25").Value = "'85.37"
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "'7.22"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'8.19"
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("D32").Value = "'28.86"
$ws.Range("E32").Value = "  +5.25%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").Value = "'6.06"
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("E37").Value = "  +5.92%  "
$ws.Range("D38").Value = "'52.35"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("D40").Value = "'432.79"
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("D41").Value = "3.115.39"
$ws.Range("E41").Value = "  +4.97%  "
$ws.Range("E42").Value = "  +8.82%  "
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("D46").Value = "'2.23"
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("D47").Value = "'36.82"
$ws.Range("E47").Value = "  +9.57%  "
$ws.Range("D48").Value = "'26.41"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "'125.97"
$ws.Range("E51").Value = "  +3.50%  "
